$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new data point was recorded as the very first sample (new row 2), pushing
# all previously recorded samples (rows 2-21) down by one row (to rows 3-22).
# Shift the existing values down by copying from the bottom row upward so we
# never overwrite data before it has been read, and avoid Insert() (which
# would copy cell formatting from the row above).
for ($r = 21; $r -ge 2; $r--) {
    $a = $ws.Cells.Item($r, 1).Value()
    $b = $ws.Cells.Item($r, 2).Value()
    $c = $ws.Cells.Item($r, 3).Value()
    $ws.Cells.Item($r + 1, 1).Value = $a
    $ws.Cells.Item($r + 1, 2).Value = $b
    $ws.Cells.Item($r + 1, 3).Value = $c
}

# Write the new first sample into row 2
$ws.Cells.Item(2, 1).Value = 0.3229818344116211
$ws.Cells.Item(2, 2).Value = 0.6911778450012207
$ws.Cells.Item(2, 3).Value = 0.0410229265689849

# Append 9 additional newly recorded samples after the existing data
# (which now ends at row 22)
$newRows = @(
    @(0.9749262332916433, 1.300361778587099, -6.267426431179062),
    @(-3.844243764877326, 1.833226948976521, -1.409952521324157),
    @(0.4797788858413697, -0.523662269115448, -1.702465817332268),
    @(1.155098915100098, 1.092013478279114, 1.727226853370667),
    @(-1.098365545272828, -0.6193101108074199, 0.1845241859555233),
    @(-0.8518145084381094, -0.03355145454406605, 0.7549576908350003),
    @(0.405293345451355, 0.8384262472391129, 0.3231545425951481),
    @(0.2438197135925255, 0.4860433936119046, -0.09267929568886754),
    @(-0.07322704792022328, 0.1344193816185026, -0.148086081258953)
)

$startRow = 23
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $vals = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $vals[0]
    $ws.Cells.Item($r, 2).Value = $vals[1]
    $ws.Cells.Item($r, 3).Value = $vals[2]
}
